$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# EditRecipient (sheet2): C2 "91827" -> "918279", selection C8 -> C2
# ---------------------------------------------------------------------------
$wsEdit = $wb.Worksheets.Item("EditRecipient")
$wsEdit.Range("C2").Value = "'918279"

# ---------------------------------------------------------------------------
# DeleteRecipient (sheet3): C2 "91827" -> "918279", selection A7 -> C2
# ---------------------------------------------------------------------------
$wsDelete = $wb.Worksheets.Item("DeleteRecipient")
$wsDelete.Range("C2").Value = "'918279"

# ---------------------------------------------------------------------------
# AddressCreate (sheet4): update B2, append rows 3 & 4, selection E6 -> C11
# New shared strings must land in this order so they match the target
# workbook: "10", "News", "82", "1000"
# ---------------------------------------------------------------------------
$wsAddr = $wb.Worksheets.Item("AddressCreate")

$wsAddr.Range("B3").Value = "'10"
$wsAddr.Range("A3").Value = "News"
$wsAddr.Range("B4").Value = "'82"
$wsAddr.Range("B2").Value = "'1000"

$wsAddr.Range("C3").Value = "'TrialData Recipient<91827>"
$wsAddr.Range("E3").Value = "'TrialData Recipient<991827>,Palak Garg<9917186286>"
$wsAddr.Range("D3").Value = "'Palak Garg<9917186286>"
$wsAddr.Range("D2").Copy()
$wsAddr.Range("D3").PasteSpecial(-4122)

$wsAddr.Range("A4").Value = "Fax Address"
$wsAddr.Range("C4").Value = "'TrialData Recipient<91827>"
$wsAddr.Range("E4").Value = "'TrialData Recipient<991827>,Palak Garg<9917186286>"
$wsAddr.Range("D4").Value = "'Palak Garg<9917186286>"
$wsAddr.Range("D2").Copy()
$wsAddr.Range("D4").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Selections (also moves the active cell on each sheet)
# ---------------------------------------------------------------------------
$wsCreate = $wb.Worksheets.Item("CreateRecipient")
$wsCreate.Activate()
$wsCreate.Range("C2").Select()

$wsEdit.Activate()
$wsEdit.Range("C2").Select()

$wsDelete.Activate()
$wsDelete.Range("C2").Select()

$wsAddr.Activate()
$wsAddr.Range("C11").Select()

# Restore the originally active sheet/tab
$wsCreate.Activate()
